# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计" (before "2022-Q2").
# 2. Populate it with the new quarter's fund-holding table (same layout/style
#    as the existing "2022-Q2" sheet).
# 3. Update the "总计" (summary) sheet: add a new top row for 2022-Q3 and
#    shift the existing rows down.
# The other quarterly sheets ("2022-Q2", "2021-Q1", "2020-Q4") keep their
# data untouched - they just shift right by one tab position automatically.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Step 1: insert the new "2022-Q3" sheet right after "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q3"

# Re-fetch the "2022-Q2" sheet by name now that tab positions have shifted.
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Copy the header-row / index-column formatting (bold + border + centered)
# from the existing "2022-Q2" sheet so the new sheet reuses the same style.
$q2Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$q2Sheet.Range("A2:A6").Copy()
$newSheet.Range("A2:A6").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Step 2: populate "2022-Q3" header + data.
# ---------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns D, E, F, G are stored as text in the source data (e.g. "5.31"),
# so force Text format before assigning to avoid numeric auto-conversion.
$newSheet.Range("D2:G6").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "900008"
$newSheet.Range("C2").Value = "中信证券稳健回报混合A"
$newSheet.Range("D2").Value = "3.87"
$newSheet.Range("E2").Value = "54.97"
$newSheet.Range("F2").Value = "6.28"
$newSheet.Range("G2").Value = "0.2430"
$newSheet.Range("H2").Value = 2

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "900078"
$newSheet.Range("C3").Value = "中信证券稳健回报混合C"
$newSheet.Range("D3").Value = "1.44"
$newSheet.Range("E3").Value = "54.97"
$newSheet.Range("F3").Value = "6.28"
$newSheet.Range("G3").Value = "0.0904"
$newSheet.Range("H3").Value = 2

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "900027"
$newSheet.Range("C4").Value = "中信证券信远一年混合A"
$newSheet.Range("D4").Value = "0.66"
$newSheet.Range("E4").Value = "63.33"
$newSheet.Range("F4").Value = "6.68"
$newSheet.Range("G4").Value = "0.0441"
$newSheet.Range("H4").Value = 2

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "900077"
$newSheet.Range("C5").Value = "中信证券信远一年混合B"
$newSheet.Range("D5").Value = "0.01"
$newSheet.Range("E5").Value = "63.33"
$newSheet.Range("F5").Value = "6.68"
$newSheet.Range("G5").Value = "0.0007"
$newSheet.Range("H5").Value = 2

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "900087"
$newSheet.Range("C6").Value = "中信证券信远一年混合C"
$newSheet.Range("D6").Value = "0.01"
$newSheet.Range("E6").Value = "63.33"
$newSheet.Range("F6").Value = "6.68"
$newSheet.Range("G6").Value = "0.0007"
$newSheet.Range("H6").Value = 2

# ---------------------------------------------------------------------
# Step 3: update the "总计" summary sheet - add the new 2022-Q3 row on
# top and push the existing rows down by one.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Give the new A5 index cell the same style as the existing index cells
# (bold + border + centered) before writing its value.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A5").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2020-Q4"
$totalSheet.Range("C5").Value = 2
$totalSheet.Range("D5").Value = 0.14

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q1"
$totalSheet.Range("C4").Value = 3
$totalSheet.Range("D4").Value = 0.98

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 5
$totalSheet.Range("D3").Value = 0.68

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 5
$totalSheet.Range("D2").Value = 0.38

Write-Host "2022-Q3 sheet added and 总计 summary updated"
